$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric (and date-serial) cell updates
$ws.Range("B2").Value = 44019
$ws.Range("C2").Value = 23607
$ws.Range("D2").Value = 689
$ws.Range("E2").Value = 2747
$ws.Range("F2").Value = 89
$ws.Range("B3").Value = 44019
$ws.Range("C3").Value = 15880
$ws.Range("D3").Value = 137
$ws.Range("B4").Value = 44019
$ws.Range("E4").Value = 33277
$ws.Range("F4").Value = 5206
$ws.Range("H4").Value = 30.51
$ws.Range("K4").Value = 110332
$ws.Range("L4").Value = 17065
$ws.Range("B6").Value = 44019
$ws.Range("C6").Value = 53514
$ws.Range("D6").Value = 665
$ws.Range("E6").Value = 10872
$ws.Range("F6").Value = 234
$ws.Range("G6").Value = 20.32
$ws.Range("H6").Value = 35.19
$ws.Range("B7").Value = 44019
$ws.Range("B8").Value = 44019
$ws.Range("C8").Value = 17519
$ws.Range("D8").Value = 602
$ws.Range("E8").Value = 1675
$ws.Range("F8").Value = 87
$ws.Range("G8").Value = 13.64
$ws.Range("H8").Value = 15.43
$ws.Range("K8").Value = 12280
$ws.Range("L8").Value = 564
$ws.Range("B9").Value = 44019
$ws.Range("B10").Value = 44019
$ws.Range("C10").Value = 17578
$ws.Range("D10").Value = 399
$ws.Range("E10").Value = 635
$ws.Range("G10").Value = 4.57
$ws.Range("H10").Value = 3.32
$ws.Range("K10").Value = 13892
$ws.Range("L10").Value = 392
$ws.Range("B12").Value = 44019
$ws.Range("C12").Value = 13727
$ws.Range("D12").Value = 519
$ws.Range("E12").Value = 261
$ws.Range("G12").Value = 1.9
$ws.Range("B13").Value = 44019
$ws.Range("C13").Value = 51058
$ws.Range("D13").Value = 1057
$ws.Range("E13").Value = 6050
$ws.Range("G13").Value = 20.55
$ws.Range("K13").Value = 29436
$ws.Range("B14").Value = 44019
$ws.Range("C14").Value = 14768
$ws.Range("E14").Value = 1902
$ws.Range("G14").Value = 24.07
$ws.Range("K14").Value = 7901
$ws.Range("B15").Value = 44019
$ws.Range("C15").Value = 45263
$ws.Range("D15").Value = 1007
$ws.Range("E15").Value = 15769
$ws.Range("F15").Value = 451
$ws.Range("G15").Value = 46.33
$ws.Range("H15").Value = 46.35
$ws.Range("K15").Value = 34039
$ws.Range("L15").Value = 973
$ws.Range("B16").Value = 44018
$ws.Range("C16").Value = 120539
$ws.Range("D16").Value = 3579
$ws.Range("E16").Value = 3232
$ws.Range("F16").Value = 369
$ws.Range("H16").Value = 11.09
$ws.Range("K16").Value = 68774
$ws.Range("L16").Value = 3327
$ws.Range("B18").Value = 44018
$ws.Range("C18").Value = 32214
$ws.Range("D18").Value = 1158
$ws.Range("E18").Value = 15391
$ws.Range("F18").Value = 579
$ws.Range("G18").Value = 47.78
$ws.Range("H18").Value = 50
$ws.Range("B19").Value = 44019
$ws.Range("C19").Value = 88691
$ws.Range("D19").Value = 6787
$ws.Range("E19").Value = 11754
$ws.Range("G19").Value = 30.08
$ws.Range("K19").Value = 39071
$ws.Range("B20").Value = 44019
$ws.Range("C20").Value = 210594
$ws.Range("D20").Value = 3841
$ws.Range("E20").Value = 27447
$ws.Range("F20").Value = 753
$ws.Range("G20").Value = 13.03
$ws.Range("H20").Value = 19.6
$ws.Range("B21").Value = 44019
$ws.Range("C21").Value = 1327
$ws.Range("G21").Value = 0.45
$ws.Range("B22").Value = 44019
$ws.Range("C22").Value = 1254
$ws.Range("E22").Value = 128
$ws.Range("G22").Value = 10.49
$ws.Range("K22").Value = 1220
$ws.Range("B23").Value = 44019
$ws.Range("C23").Value = 34664
$ws.Range("D23").Value = 1696
$ws.Range("E23").Value = 1827
$ws.Range("H23").Value = 6.73
$ws.Range("K23").Value = 28456
$ws.Range("L23").Value = 1634
$ws.Range("B24").Value = 44019
$ws.Range("C24").Value = 20201
$ws.Range("D24").Value = 282
$ws.Range("E24").Value = 1200
$ws.Range("G24").Value = 7.67
$ws.Range("K24").Value = 15643
$ws.Range("B25").Value = 44019
$ws.Range("C25").Value = 66540
$ws.Range("D25").Value = 5926
$ws.Range("E25").Value = 20100
$ws.Range("F25").Value = 2364
$ws.Range("G25").Value = 30.21
$ws.Range("H25").Value = 39.89
$ws.Range("B26").Value = 44018
$ws.Range("C26").Value = 277774
$ws.Range("D26").Value = 6361
$ws.Range("E26").Value = 7869
$ws.Range("F26").Value = 577
$ws.Range("G26").Value = 4.4
$ws.Range("K26").Value = 180178
$ws.Range("L26").Value = 6282
$ws.Range("B27").Value = 44019
$ws.Range("C27").Value = 48626
$ws.Range("D27").Value = 2524
$ws.Range("E27").Value = 5883
$ws.Range("F27").Value = 367
$ws.Range("G27").Value = 12.1
$ws.Range("H27").Value = 14.54
$ws.Range("B28").Value = 44019
$ws.Range("C28").Value = 1184
$ws.Range("D28").Value = 17
$ws.Range("E28").Value = 29
$ws.Range("G28").Value = 1.4
$ws.Range("K28").Value = 2069
$ws.Range("L28").Value = 34
$ws.Range("B29").Value = 44019
$ws.Range("C29").Value = 32556
$ws.Range("D29").Value = 805
$ws.Range("E29").Value = 5560
$ws.Range("F29").Value = 192
$ws.Range("G29").Value = 18.97
$ws.Range("H29").Value = 24.21
$ws.Range("K29").Value = 29315
$ws.Range("L29").Value = 793
$ws.Range("B30").Value = 44019
$ws.Range("C30").Value = 100470
$ws.Range("D30").Value = 2899
$ws.Range("E30").Value = 27660
$ws.Range("F30").Value = 1361
$ws.Range("G30").Value = 27.53
$ws.Range("H30").Value = 46.95
$ws.Range("B31").Value = 44019
$ws.Range("C31").Value = 37420
$ws.Range("D31").Value = 1384
$ws.Range("E31").Value = 1464
$ws.Range("G31").Value = 5.48
$ws.Range("H31").Value = 3.43
$ws.Range("K31").Value = 26717
$ws.Range("L31").Value = 1312
$ws.Range("B33").Value = 44019
$ws.Range("C33").Value = 12414
$ws.Range("D33").Value = 514
$ws.Range("E33").Value = 3174
$ws.Range("G33").Value = 25.57
$ws.Range("H33").Value = 25.49
$ws.Range("B35").Value = 44019
$ws.Range("C35").Value = 75875
$ws.Range("D35").Value = 1420
$ws.Range("E35").Value = 12230
$ws.Range("F35").Value = 451
$ws.Range("G35").Value = 23.82
$ws.Range("H35").Value = 32.92
$ws.Range("K35").Value = 51350
$ws.Range("L35").Value = 1370
$ws.Range("B36").Value = 44019
$ws.Range("C36").Value = 148452
$ws.Range("D36").Value = 7063
$ws.Range("E36").Value = 24901
$ws.Range("F36").Value = 1965
$ws.Range("G36").Value = 16.77
$ws.Range("H36").Value = 27.82
$ws.Range("B37").Value = 44019
$ws.Range("C37").Value = 8539
$ws.Range("E37").Value = 127
$ws.Range("B38").Value = 44019
$ws.Range("C38").Value = 39133
$ws.Range("D38").Value = 1477
$ws.Range("E38").Value = 8000
$ws.Range("F38").Value = 130
$ws.Range("G38").Value = 20.44
$ws.Range("H38").Value = 8.800000000000001
$ws.Range("B39").Value = 44019
$ws.Range("C39").Value = 110338
$ws.Range("D39").Value = 8213
$ws.Range("E39").Value = 10400
$ws.Range("F39").Value = 672
$ws.Range("G39").Value = 9.43
$ws.Range("B40").Value = 44019
$ws.Range("C40").Value = 12577
$ws.Range("D40").Value = 358
$ws.Range("E40").Value = 3635
$ws.Range("F40").Value = 145
$ws.Range("G40").Value = 31.14
$ws.Range("H40").Value = 40.5
$ws.Range("K40").Value = 11673
$ws.Range("L40").Value = 358
$ws.Range("B41").Value = 44019
$ws.Range("C41").Value = 24629
$ws.Range("D41").Value = 1042
$ws.Range("C42").Value = 32042
$ws.Range("D42").Value = 728
$ws.Range("E42").Value = 2882
$ws.Range("H42").Value = 4.95

# String-valued cells (numeric-looking text must stay text)
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "214371"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18611"
$ws.Range("D4").ClearFormats()
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "26033"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "194"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "666"
$ws.Range("E7").ClearFormats()

Write-Output "applied covid disparities update for 2020-07-07 run"
